# dictionary_tree.xlsx -- "Adjust indentation and alignment"
#
# The sheet has a small list of folders (column D) with a one-line
# description next to each (column E). The GBDT/random-forest row's
# description had an extra leading space before "GBDT" (" >>>  GBDT, ..."
# with two spaces) -- fix it to a single space (" >>> GBDT, ...").
# Re-writing the cell text naturally moves its shared-string entry to the
# end of the table and shifts every other description's shared-string
# index down by one, which is exactly what the rest of the diff shows.
#
# The author's cursor also ended up resting on B19 instead of E14 when the
# file was saved, so the sheet's stored selection is updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = " >>> GBDT, random forest etc. "

$ws.Range("B19").Select()
